$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.759.67'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '2.452.32'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.06'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.91'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.07%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.505'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '2.453.25'
$ws.Range('E9').Value = '  -1.20%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.150'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.67%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.333'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.63%  '
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = '2.905.82'
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('D15').Value = '68.607.87'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000168'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '23.55'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').Value = '2.447.48'
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.71'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '337.51'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.02'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.78'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.92'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.77'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.29%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.68'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.36%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.582.10'
$ws.Range('E27').Value = '  -1.35%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.52%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.18'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.41%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0815'
$ws.Range('E30').Value = '  -5.74%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.17'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.11%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '433.35'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.14'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.76%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.61'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.50%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '156.52'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.68%  '
$ws.Range('B37').Value = 'WhiteBITCoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.03'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('B38').Value = 'USDe'
$ws.Range('C38').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.108'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.76%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.76'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.81%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.300'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.41'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.65%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '37.38'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.47'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -6.36%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.10'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.28%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.06'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -4.31%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '132.42'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.48%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.33'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.53%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0715'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.03%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.481'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.23%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.558'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.13%  '
